$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "baseline" row to describe the baseline SVM model and its parameters
$ws.Range("A2").Value = "baseline(svm)"
$ws.Range("B2").Value = "C=1, gamma=10, kernel='rbf'"

# Move the active selection as in the authored workbook
$ws.Range("F2").Select()
